$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 907 (shifts existing rows 907:1010 down to 908:1011)
$ws.Rows.Item(907).Insert()

# Populate the newly inserted row 907 with the new data record
$ws.Cells.Item(907, 1).Value = 3
$ws.Cells.Item(907, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(907, 3).Value = "Coquimbo"
$ws.Cells.Item(907, 4).Value = 45212
$ws.Cells.Item(907, 5).Value = 5
$ws.Cells.Item(907, 6).Value = 100112045
$ws.Cells.Item(907, 7).Value = "Zapallo"
$ws.Cells.Item(907, 8).Value = "Camote"
$ws.Cells.Item(907, 9).Value = "1a (guarda)"
$ws.Cells.Item(907, 10).Value = 120
$ws.Cells.Item(907, 11).Value = 1000
$ws.Cells.Item(907, 12).Value = 1000
$ws.Cells.Item(907, 13).Value = 1000
$ws.Cells.Item(907, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(907, 15).Value = "Provincia de Talca"
$ws.Cells.Item(907, 16).Value = 1000
$ws.Cells.Item(907, 17).Value = 1
$ws.Cells.Item(907, 18).Value = "Hortaliza"
